# Add the "countdown track" (Index.xlsx) as a new entry in the file list.
# This inserts a new row at position 27 (pushing the remaining rows down by
# one, so the former row 73 becomes row 74) and then appends the ".mp3"
# extension to every pre-existing filename in column A.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 27, shifting rows 27-73 down to 28-74.
$ws.Rows(27).Insert()

# Populate the newly inserted row with the new entry.
$ws.Cells.Item(27, 1).Value2 = "Index.xlsx"

# Make sure the rest of the row exists as blank cells, matching the layout
# of every other data row (Artiste/Song/In/Out left blank).
$ws.Cells.Item(27, 2).Style = "Normal"
$ws.Cells.Item(27, 3).Style = "Normal"
$ws.Cells.Item(27, 4).Style = "Normal"
$ws.Cells.Item(27, 5).Style = "Normal"

# Append the ".mp3" extension to every other filename in column A
# (rows 2-26 keep their original row number, rows 28-74 are the old
# rows 27-73 that were shifted down by the insert above).
for ($r = 2; $r -le 74; $r++) {
    if ($r -eq 27) { continue }
    $cell = $ws.Cells.Item($r, 1)
    $name = $cell.Value2
    $cell.Value2 = $name + ".mp3"
}
